# Appends a new data row (row 23) mirroring the existing Adafruit IO feed
# rows, extending the sheet's used range from A1:F22 to A1:F23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 23

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"

# "25" looks numeric, so Excel would otherwise coerce it to a Number cell.
# The source data stores it as text (matching the other rows in column C),
# so force a text number format before assigning, then restore the Normal
# style so no stray formatting is left behind on the cell.
$ws.Range("C$newRow").NumberFormat = "@"
$ws.Range("C$newRow").Value = "25"
$ws.Range("C$newRow").Style = "Normal"

$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
